$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (shifts rows 4-16 down to 5-17, and
# auto-extends the SUM formula range used in the total row).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with the new timesheet entry.
$ws.Range("A4").Value2 = 44632
$ws.Range("B4").Value2 = 0.5

# A few of the existing dates further down were corrected.
$ws.Range("A8").Value2 = 44644
$ws.Range("A14").Value2 = 44660
$ws.Range("A15").Value2 = 44661
$ws.Range("A16").Value2 = 44662

# Update the selection to match the new active cell.
$ws.Range("E17").Select() | Out-Null
